$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Email column (D) for rows 2-11 with new addresses ---
$ws.Range("D2").Value = "johnsmith99@gmail.com"
$ws.Range("D3").Value = "amyjackson@gmail.com"
$ws.Range("D4").Value = "jadewilliam@gmail.com"
$ws.Range("D5").Value = "cathiecat@gmail.com"
$ws.Range("D6").Value = "mysonjacob@gmail.com"
$ws.Range("D7").Value = "willy9898@gmail.com"
$ws.Range("D8").Value = "jadejady@gmail.com"
$ws.Range("D9").Value = "smartemily@gmail.com"
$ws.Range("D10").Value = "naughtynatasha@gmail.com"
$ws.Range("D11").Value = "danieldanny77@gmail.com"

# --- Rebuild the mailto hyperlinks so they point at the new addresses ---
# (the engine only supports clearing the whole collection at once, so
#  delete everything and re-add all 14 links, keeping the 4 untouched
#  ones - D12,D13,D14,D15 - first, then the 10 refreshed ones in sheet
#  order to match the target layout)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:Dali@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:Daley@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:Gary@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:emily@gmail.com")

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:johnsmith99@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:amyjackson@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jadewilliam@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:cathiecat@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:mysonjacob@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:willy9898@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:jadejady@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:smartemily@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:danieldanny77@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:naughtynatasha@gmail.com")

# --- Widen column D to fit the longer email addresses ---
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668

# --- Move the active selection to C16 ---
$ws.Range("C16").Select()
